$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.101967
$ws.Range("H2").Value = 27.305901
$ws.Range("I2").Value = 0.2872601673725235
$ws.Range("J2").Value = 0.2872601673725235
$ws.Range("M2").Value = 3.889188333333333
$ws.Range("N2").Value = 11.667565
$ws.Range("O2").Value = 0.04320782608967438
$ws.Range("P2").Value = 0.04320782608967438
$ws.Range("Q2").Value = 35.399263866785
$ws.Range("R2").Value = 318.593374801065
$ws.Range("S2").Value = 0.01241188735432275
$ws.Range("T2").Value = 0.01241188735432275

$ws.Range("G3").Value = 9.101967
$ws.Range("H3").Value = 27.305901
$ws.Range("I3").Value = 0.2872601673725235
$ws.Range("J3").Value = 0.2872601673725235
$ws.Range("O3").Value = 0.86451478461177
$ws.Range("P3").Value = 0.8645147846117699
$ws.Range("Q3").Value = 708.2787945335269
$ws.Range("R3").Value = 6374.509150801742
$ws.Range("S3").Value = 0.2483406617235981
$ws.Range("T3").Value = 0.2483406617235981

$ws.Range("G4").Value = 9.101967
$ws.Range("H4").Value = 27.305901
$ws.Range("I4").Value = 0.2872601673725235
$ws.Range("J4").Value = 0.2872601673725235
$ws.Range("M4").Value = 0.359731
$ws.Range("N4").Value = 1.079193
$ws.Range("O4").Value = 0.003996513707975397
$ws.Range("P4").Value = 0.003996513707975397
$ws.Range("Q4").Value = 3.274259690877
$ws.Range("R4").Value = 29.468337217893
$ws.Range("S4").Value = 0.001148039196659597
$ws.Range("T4").Value = 0.001148039196659597

$ws.Range("G5").Value = 9.101967
$ws.Range("H5").Value = 27.305901
$ws.Range("I5").Value = 0.2872601673725235
$ws.Range("J5").Value = 0.2872601673725235
$ws.Range("M5").Value = 7.581867
$ws.Range("N5").Value = 22.745601
$ws.Range("O5").Value = 0.08423248315420773
$ws.Range("P5").Value = 0.08423248315420773
$ws.Range("Q5").Value = 69.00990323238901
$ws.Range("R5").Value = 621.089129091501
$ws.Range("S5").Value = 0.02419663720908097
$ws.Range("T5").Value = 0.02419663720908097

$ws.Range("G6").Value = 9.101967
$ws.Range("H6").Value = 27.305901
$ws.Range("I6").Value = 0.2872601673725235
$ws.Range("J6").Value = 0.2872601673725235
$ws.Range("M6").Value = 0.3644006666666667
$ws.Range("N6").Value = 1.093202
$ws.Range("O6").Value = 0.004048392436372474
$ws.Range("P6").Value = 0.004048392436372474
$ws.Range("Q6").Value = 3.316762842778
$ws.Range("R6").Value = 29.850865585002
$ws.Range("S6").Value = 0.001162941888862015
$ws.Range("T6").Value = 0.001162941888862015

$ws.Range("I7").Value = 0.3055950511371977
$ws.Range("J7").Value = 0.3055950511371977
$ws.Range("M7").Value = 3.889188333333333
$ws.Range("N7").Value = 11.667565
$ws.Range("O7").Value = 0.04320782608967438
$ws.Range("P7").Value = 0.04320782608967438
$ws.Range("Q7").Value = 37.65868393984667
$ws.Range("R7").Value = 338.92815545862
$ws.Range("S7").Value = 0.01320409782340119
$ws.Range("T7").Value = 0.01320409782340119

$ws.Range("I8").Value = 0.3055950511371977
$ws.Range("J8").Value = 0.3055950511371977
$ws.Range("O8").Value = 0.86451478461177
$ws.Range("P8").Value = 0.8645147846117699
$ws.Range("S8").Value = 0.2641914398122973
$ws.Range("T8").Value = 0.2641914398122973

$ws.Range("I9").Value = 0.3055950511371977
$ws.Range("J9").Value = 0.3055950511371977
$ws.Range("M9").Value = 0.359731
$ws.Range("N9").Value = 1.079193
$ws.Range("O9").Value = 0.003996513707975397
$ws.Range("P9").Value = 0.003996513707975397
$ws.Range("Q9").Value = 3.483245055596
$ws.Range("R9").Value = 31.349205500364
$ws.Range("S9").Value = 0.001221314810959253
$ws.Range("T9").Value = 0.001221314810959253

$ws.Range("I10").Value = 0.3055950511371977
$ws.Range("J10").Value = 0.3055950511371977
$ws.Range("M10").Value = 7.581867
$ws.Range("N10").Value = 22.745601
$ws.Range("O10").Value = 0.08423248315420773
$ws.Range("P10").Value = 0.08423248315420773
$ws.Range("Q10").Value = 73.414581284172
$ws.Range("R10").Value = 660.7312315575481
$ws.Range("S10").Value = 0.02574102999692325
$ws.Range("T10").Value = 0.02574102999692325

$ws.Range("I11").Value = 0.3055950511371977
$ws.Range("J11").Value = 0.3055950511371977
$ws.Range("M11").Value = 0.3644006666666667
$ws.Range("N11").Value = 1.093202
$ws.Range("O11").Value = 0.004048392436372474
$ws.Range("P11").Value = 0.004048392436372474
$ws.Range("Q11").Value = 3.528461045677334
$ws.Range("R11").Value = 31.756149411096
$ws.Range("S11").Value = 0.001237168693616691
$ws.Range("T11").Value = 0.001237168693616691

$ws.Range("G12").Value = 3.905093666666666
$ws.Range("H12").Value = 11.715281
$ws.Range("I12").Value = 0.1232456523180152
$ws.Range("J12").Value = 0.1232456523180152
$ws.Range("M12").Value = 3.889188333333333
$ws.Range("N12").Value = 11.667565
$ws.Range("O12").Value = 0.04320782608967438
$ws.Range("P12").Value = 0.04320782608967438
$ws.Range("Q12").Value = 15.18764472897389
$ws.Range("R12").Value = 136.688802560765
$ws.Range("S12").Value = 0.005325176711665274
$ws.Range("T12").Value = 0.005325176711665274

$ws.Range("G13").Value = 3.905093666666666
$ws.Range("H13").Value = 11.715281
$ws.Range("I13").Value = 0.1232456523180152
$ws.Range("J13").Value = 0.1232456523180152
$ws.Range("O13").Value = 0.86451478461177
$ws.Range("P13").Value = 0.8645147846117699
$ws.Range("Q13").Value = 303.8788247383425
$ws.Range("R13").Value = 2734.909422645082
$ws.Range("S13").Value = 0.106547688568046
$ws.Range("T13").Value = 0.106547688568046

$ws.Range("G14").Value = 3.905093666666666
$ws.Range("H14").Value = 11.715281
$ws.Range("I14").Value = 0.1232456523180152
$ws.Range("J14").Value = 0.1232456523180152
$ws.Range("M14").Value = 0.359731
$ws.Range("N14").Value = 1.079193
$ws.Range("O14").Value = 0.003996513707975397
$ws.Range("P14").Value = 0.003996513707975397
$ws.Range("Q14").Value = 1.404783249803667
$ws.Range("R14").Value = 12.643049248233
$ws.Range("S14").Value = 0.0004925529389373176
$ws.Range("T14").Value = 0.0004925529389373176

$ws.Range("G15").Value = 3.905093666666666
$ws.Range("H15").Value = 11.715281
$ws.Range("I15").Value = 0.1232456523180152
$ws.Range("J15").Value = 0.1232456523180152
$ws.Range("M15").Value = 7.581867
$ws.Range("N15").Value = 22.745601
$ws.Range("O15").Value = 0.08423248315420773
$ws.Range("P15").Value = 0.08423248315420773
$ws.Range("Q15").Value = 29.607900803209
$ws.Range("R15").Value = 266.471107228881
$ws.Range("S15").Value = 0.01038128733270656
$ws.Range("T15").Value = 0.01038128733270656

$ws.Range("G16").Value = 3.905093666666666
$ws.Range("H16").Value = 11.715281
$ws.Range("I16").Value = 0.1232456523180152
$ws.Range("J16").Value = 0.1232456523180152
$ws.Range("M16").Value = 0.3644006666666667
$ws.Range("N16").Value = 1.093202
$ws.Range("O16").Value = 0.004048392436372474
$ws.Range("P16").Value = 0.004048392436372474
$ws.Range("Q16").Value = 1.423018735529111
$ws.Range("R16").Value = 12.807168619762
$ws.Range("S16").Value = 0.0004989467666600445
$ws.Range("T16").Value = 0.0004989467666600445

$ws.Range("G17").Value = 6.285238333333333
$ws.Range("H17").Value = 18.855715
$ws.Range("I17").Value = 0.1983635642284282
$ws.Range("J17").Value = 0.1983635642284282
$ws.Range("M17").Value = 3.889188333333333
$ws.Range("N17").Value = 11.667565
$ws.Range("O17").Value = 0.04320782608967438
$ws.Range("P17").Value = 0.04320782608967438
$ws.Range("Q17").Value = 24.44447559821944
$ws.Range("R17").Value = 220.000280383975
$ws.Range("S17").Value = 0.008570858385709878
$ws.Range("T17").Value = 0.008570858385709878

$ws.Range("G18").Value = 6.285238333333333
$ws.Range("H18").Value = 18.855715
$ws.Range("I18").Value = 0.1983635642284282
$ws.Range("J18").Value = 0.1983635642284282
$ws.Range("O18").Value = 0.86451478461177
$ws.Range("P18").Value = 0.8645147846117699
$ws.Range("Q18").Value = 489.0921962350827
$ws.Range("R18").Value = 4401.829766115745
$ws.Range("S18").Value = 0.1714882340037626
$ws.Range("T18").Value = 0.1714882340037626

$ws.Range("G19").Value = 6.285238333333333
$ws.Range("H19").Value = 18.855715
$ws.Range("I19").Value = 0.1983635642284282
$ws.Range("J19").Value = 0.1983635642284282
$ws.Range("M19").Value = 0.359731
$ws.Range("N19").Value = 1.079193
$ws.Range("O19").Value = 0.003996513707975397
$ws.Range("P19").Value = 0.003996513707975397
$ws.Range("Q19").Value = 2.260995070888333
$ws.Range("R19").Value = 20.348955637995
$ws.Range("S19").Value = 0.0007927627036017714
$ws.Range("T19").Value = 0.0007927627036017714

$ws.Range("G20").Value = 6.285238333333333
$ws.Range("H20").Value = 18.855715
$ws.Range("I20").Value = 0.1983635642284282
$ws.Range("J20").Value = 0.1983635642284282
$ws.Range("M20").Value = 7.581867
$ws.Range("N20").Value = 22.745601
$ws.Range("O20").Value = 0.08423248315420773
$ws.Range("P20").Value = 0.08423248315420773
$ws.Range("Q20").Value = 47.653841106635
$ws.Range("R20").Value = 428.884569959715
$ws.Range("S20").Value = 0.01670865558227968
$ws.Range("T20").Value = 0.01670865558227968

$ws.Range("G21").Value = 6.285238333333333
$ws.Range("H21").Value = 18.855715
$ws.Range("I21").Value = 0.1983635642284282
$ws.Range("J21").Value = 0.1983635642284282
$ws.Range("M21").Value = 0.3644006666666667
$ws.Range("N21").Value = 1.093202
$ws.Range("O21").Value = 0.004048392436372474
$ws.Range("P21").Value = 0.004048392436372474
$ws.Range("Q21").Value = 2.290345038825556
$ws.Range("R21").Value = 20.61310534943
$ws.Range("S21").Value = 0.0008030535530742541
$ws.Range("T21").Value = 0.0008030535530742541

$ws.Range("G22").Value = 2.710232666666667
$ws.Range("H22").Value = 8.130698000000001
$ws.Range("I22").Value = 0.08553556494383548
$ws.Range("J22").Value = 0.08553556494383548
$ws.Range("M22").Value = 3.889188333333333
$ws.Range("N22").Value = 11.667565
$ws.Range("O22").Value = 0.04320782608967438
$ws.Range("P22").Value = 0.04320782608967438
$ws.Range("Q22").Value = 10.54060526781889
$ws.Range("R22").Value = 94.86544741037001
$ws.Range("S22").Value = 0.003695805814575292
$ws.Range("T22").Value = 0.003695805814575292

$ws.Range("G23").Value = 2.710232666666667
$ws.Range("H23").Value = 8.130698000000001
$ws.Range("I23").Value = 0.08553556494383548
$ws.Range("J23").Value = 0.08553556494383548
$ws.Range("O23").Value = 0.86451478461177
$ws.Range("P23").Value = 0.8645147846117699
$ws.Range("Q23").Value = 210.8995040360016
$ws.Range("R23").Value = 1898.095536324014
$ws.Range("S23").Value = 0.073946760504066
$ws.Range("T23").Value = 0.07394676050406598

$ws.Range("G24").Value = 2.710232666666667
$ws.Range("H24").Value = 8.130698000000001
$ws.Range("I24").Value = 0.08553556494383548
$ws.Range("J24").Value = 0.08553556494383548
$ws.Range("M24").Value = 0.359731
$ws.Range("N24").Value = 1.079193
$ws.Range("O24").Value = 0.003996513707975397
$ws.Range("P24").Value = 0.003996513707975397
$ws.Range("Q24").Value = 0.9749547074126668
$ws.Range("R24").Value = 8.774592366714002
$ws.Range("S24").Value = 0.0003418440578174583
$ws.Range("T24").Value = 0.0003418440578174583

$ws.Range("G25").Value = 2.710232666666667
$ws.Range("H25").Value = 8.130698000000001
$ws.Range("I25").Value = 0.08553556494383548
$ws.Range("J25").Value = 0.08553556494383548
$ws.Range("M25").Value = 7.581867
$ws.Range("N25").Value = 22.745601
$ws.Range("O25").Value = 0.08423248315420773
$ws.Range("P25").Value = 0.08423248315420773
$ws.Range("Q25").Value = 20.548623617722
$ws.Range("R25").Value = 184.937612559498
$ws.Range("S25").Value = 0.007204873033217263
$ws.Range("T25").Value = 0.007204873033217263

$ws.Range("G26").Value = 2.710232666666667
$ws.Range("H26").Value = 8.130698000000001
$ws.Range("I26").Value = 0.08553556494383548
$ws.Range("J26").Value = 0.08553556494383548
$ws.Range("M26").Value = 0.3644006666666667
$ws.Range("N26").Value = 1.093202
$ws.Range("O26").Value = 0.004048392436372474
$ws.Range("P26").Value = 0.004048392436372474
$ws.Range("Q26").Value = 0.9876105905551111
$ws.Range("R26").Value = 8.888495314996002
$ws.Range("S26").Value = 0.0003462815341594701
$ws.Range("T26").Value = 0.0003462815341594701
